# Regenerate save_data to use K (column G) instead of Strike#.
# Updates the K column values for rows 2-15 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
